$d = $word.ActiveDocument

function Set-ParagraphXmlByNeedle($needle, $innerXml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find paragraph containing: $needle"
    }
    $rng.Expand(4) | Out-Null   ; # wdParagraph -> grow the range to the whole enclosing paragraph
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml)
}

# "(5 points) As a developer, I want to make good, consistent commits."
# -> whole paragraph highlighted yellow.
Set-ParagraphXmlByNeedle "I want to make good, consistent commits" (
    '<w:p>' +
        '<w:r><w:rPr><w:b/><w:highlight w:val="yellow"/></w:rPr><w:t>(5 points)</w:t></w:r>' +
        '<w:r><w:rPr><w:b/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>As a developer, I want to make good, consistent commits.</w:t></w:r>' +
    '</w:p>'
)

# "(5 points) As a superhero fan, I want to view a list of superheroes."
# -> highlighted, but the trailing period stays un-highlighted (its own run).
Set-ParagraphXmlByNeedle "I want to view a list of superheroes" (
    '<w:p>' +
        '<w:r><w:rPr><w:b/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">(5 points) </w:t></w:r>' +
        '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>As a superhero fan, I want to view a list of superheroes</w:t></w:r>' +
        '<w:r><w:t>.</w:t></w:r>' +
    '</w:p>'
)

# "(10 points) As a superhero fan, I want to click on a superhero name from the list and view that."
# -> no highlighting change; runs reshuffled so the bookmark sits right after "points)"
#    and the sentence + trailing period become a single run.
Set-ParagraphXmlByNeedle "I want to click on a superhero name" (
    '<w:p>' +
        '<w:r><w:rPr><w:b/></w:rPr><w:t>(</w:t></w:r>' +
        '<w:r><w:rPr><w:b/></w:rPr><w:t>10</w:t></w:r>' +
        '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> points)</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
        '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r><w:t>As a superhero fan, I want to click on a superhero name from the list and view that.</w:t></w:r>' +
    '</w:p>'
)

# "(10 points) As a superhero fan, I want to create a new superhero with ..."
# -> whole paragraph highlighted yellow.
Set-ParagraphXmlByNeedle "create a new superhero with a superhero name" (
    '<w:p>' +
        '<w:r><w:rPr><w:b/><w:highlight w:val="yellow"/></w:rPr><w:t>(</w:t></w:r>' +
        '<w:r><w:rPr><w:b/><w:highlight w:val="yellow"/></w:rPr><w:t>10</w:t></w:r>' +
        '<w:r><w:rPr><w:b/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> points) </w:t></w:r>' +
        '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">As a superhero fan, I want to </w:t></w:r>' +
        '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>create a new superhero with a superhero name, alter ego name, primary superhero ability, secondary superhero ability, and catchphrase.</w:t></w:r>' +
    '</w:p>'
)

# "(10 points) As a superhero fan, I want to edit the superhero name, ... of an existing superhero."
# -> whole paragraph highlighted yellow.
Set-ParagraphXmlByNeedle "edit the superhero name, alter ego name" (
    '<w:p>' +
        '<w:r><w:rPr><w:b/><w:highlight w:val="yellow"/></w:rPr><w:t>(</w:t></w:r>' +
        '<w:r><w:rPr><w:b/><w:highlight w:val="yellow"/></w:rPr><w:t>10</w:t></w:r>' +
        '<w:r><w:rPr><w:b/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> points) </w:t></w:r>' +
        '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">As a superhero fan, I want to edit the superhero name, alter ego name, primary superhero ability, secondary superhero ability, and catchphrase of an existing </w:t></w:r>' +
        '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>superhero.</w:t></w:r>' +
    '</w:p>'
)

# "(10 points) As a superhero fan, I want to delete a superhero from the database."
# -> whole paragraph highlighted yellow.
Set-ParagraphXmlByNeedle "I want to delete a superhero from the database" (
    '<w:p>' +
        '<w:r><w:rPr><w:b/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">(10 points) </w:t></w:r>' +
        '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>As a superhero fan, I want to delete a superhero from the database.</w:t></w:r>' +
    '</w:p>'
)
